$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.444.17"
$ws.Range("E2").Value = "  -3.03%  "

$ws.Range("D3").Value = "3.499.97"
$ws.Range("E3").Value = "  -4.72%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "606.14"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "149.19"
$ws.Range("E6").Value = "  -6.42%  "

$ws.Range("D7").Value = "3.498.53"
$ws.Range("E7").Value = "  -4.67%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -3.95%  "

$ws.Range("D11").Value = "6.94"
$ws.Range("E11").Value = "  -3.40%  "

$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -4.01%  "

$ws.Range("E13").Value = "  -4.27%  "

$ws.Range("D14").Value = "4.095.77"
$ws.Range("E14").Value = "  -4.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.85%  "

$ws.Range("D16").Value = "3.503.37"
$ws.Range("E16").Value = "  -3.95%  "

$ws.Range("D17").Value = "67.373.82"
$ws.Range("E17").Value = "  -3.20%  "

$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.52%  "

$ws.Range("D21").Value = "446.02"
$ws.Range("E21").Value = "  -5.02%  "

$ws.Range("D22").Value = "8.96"
$ws.Range("E22").Value = "  -13.02%  "

$ws.Range("D23").Value = "0.619"
$ws.Range("E23").Value = "  -4.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "

$ws.Range("E25").Value = "  +5.11%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "3.640.52"
$ws.Range("E27").Value = "  -4.67%  "

$ws.Range("D28").Value = "10.18"
$ws.Range("E28").Value = "  -7.97%  "

$ws.Range("D29").Value = "8.24"
$ws.Range("E29").Value = "  -5.52%  "

$ws.Range("E30").Value = "  -4.29%  "

$ws.Range("E31").Value = "  -7.36%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("E33").Value = "  +0.85%  "

$ws.Range("D34").Value = "25.64"
$ws.Range("E34").Value = "  -3.63%  "

$ws.Range("D35").Value = "6.14"
$ws.Range("E35").Value = "  -3.75%  "

$ws.Range("E36").Value = "  -6.48%  "

$ws.Range("D37").Value = "3.491.03"
$ws.Range("E37").Value = "  -4.95%  "

$ws.Range("D38").Value = "7.97"
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("D41").Value = "2.19"
$ws.Range("E41").Value = "  -1.06%  "

$ws.Range("D42").Value = "173.41"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").Value = "0.0875"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  -5.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.880"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.90%  "

$ws.Range("D46").Value = "45.42"
$ws.Range("E46").Value = "  -2.99%  "

$ws.Range("D47").Value = "27.88"
$ws.Range("E47").Value = "  -4.52%  "

$ws.Range("E48").Value = "  +5.60%  "

$ws.Range("D49").Value = "2.55"
$ws.Range("E49").Value = "  -5.70%  "

$ws.Range("D50").Value = "7.53"
$ws.Range("E50").Value = "  -4.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.36%  "
